# Apply updated crypto price/volume values to the worksheet.
# Column D values that look like plain numbers (e.g. "300.10") would be
# auto-converted to numeric cells by Excel, so those are written with a
# leading apostrophe (forcing text) and then have their style reset back
# to "Normal" so no stray quote-prefix style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.000.77'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '2.304.48'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''300.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '

$ws.Range("D6").Value = '''97.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").Value = '''0.512'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.62%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -2.55%  '

$ws.Range("D10").Value = '''35.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").Value = '''18.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("E13").Value = '  +1.39%  '

$ws.Range("E14").Value = '  -1.44%  '

$ws.Range("D15").Value = '2.662.29'
$ws.Range("E15").Value = '  -0.11%  '

$ws.Range("D16").Value = '2.300.11'
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("E17").Value = '  -0.93%  '

$ws.Range("D18").Value = '42.924.95'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").Value = '''12.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.62%  '

$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("E21").Value = '  -1.76%  '

$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").Value = '''236.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("E24").Value = '  -1.57%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '''2.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.57%  '

$ws.Range("D27").Value = '''4.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.34%  '

$ws.Range("D28").Value = '''25.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.74%  '

$ws.Range("D29").Value = '''165.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("E31").Value = '  -1.05%  '

$ws.Range("D32").Value = '''33.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.34%  '

$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = '''4.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '''1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("E35").Value = '  -2.77%  '

$ws.Range("D36").Value = '''16.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.20%  '

$ws.Range("E37").Value = '  -1.22%  '

$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("E39").Value = '  -1.01%  '

$ws.Range("E40").Value = '  -1.63%  '

$ws.Range("E41").Value = '  -1.58%  '

$ws.Range("E42").Value = '  -0.92%  '

$ws.Range("D43").Value = '2.008.63'
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("E45").Value = '  -0.29%  '

$ws.Range("D46").Value = '''2.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.98%  '

$ws.Range("D47").Value = '''17.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.63%  '

$ws.Range("E48").Value = '  -1.31%  '

$ws.Range("E49").Value = '  +7.31%  '

$ws.Range("D50").Value = '''53.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.52%  '

$ws.Range("D51").Value = '2.530.93'
$ws.Range("E51").Value = '  -0.01%  '
